# This script updates the "Price" (column D) and "Volume(1h)" (column E)
# values in the cryptocurrency price sheet to match the refreshed data
# pulled by the "Updated symbol list" GitHub Actions job.
#
# The sheet stores these figures as plain text (e.g. "330.80", "1.53%"),
# including values that look numeric but carry significant trailing
# zeros / percent signs that must be preserved exactly. Setting
# NumberFormat = "@" (Text) before assigning Value keeps Excel from
# re-interpreting the string as a number (which would normalize
# "330.80" -> 330.8, "1.53%" -> 0.0153, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated list of (row, column, new text value) updates for the
# "Price" (D) and "Volume(1h)" (E) columns, derived from the commit diff.
$updates = @(
    @{ Row = 2; Col = "D"; Value = "330.80" }
    @{ Row = 2; Col = "E"; Value = "1.53%" }
    @{ Row = 3; Col = "E"; Value = "0.32%" }
    @{ Row = 4; Col = "E"; Value = "-1.27%" }
    @{ Row = 5; Col = "E"; Value = "0.35%" }
    @{ Row = 6; Col = "D"; Value = "2.082" }
    @{ Row = 6; Col = "E"; Value = "10.75%" }
    @{ Row = 7; Col = "D"; Value = "2.659" }
    @{ Row = 7; Col = "E"; Value = "2.62%" }
    @{ Row = 8; Col = "D"; Value = "0.9536" }
    @{ Row = 8; Col = "E"; Value = "1.10%" }
    @{ Row = 9; Col = "D"; Value = "0.1149" }
    @{ Row = 9; Col = "E"; Value = "-0.81%" }
    @{ Row = 10; Col = "D"; Value = "0.1897" }
    @{ Row = 10; Col = "E"; Value = "3.65%" }
    @{ Row = 11; Col = "D"; Value = "10.21" }
    @{ Row = 11; Col = "E"; Value = "18.03%" }
    @{ Row = 12; Col = "D"; Value = "0.09987" }
    @{ Row = 12; Col = "E"; Value = "2.80%" }
    @{ Row = 13; Col = "D"; Value = "0.04804" }
    @{ Row = 13; Col = "E"; Value = "9.85%" }
    @{ Row = 14; Col = "E"; Value = "0.06%" }
    @{ Row = 15; Col = "D"; Value = "0.001269" }
    @{ Row = 15; Col = "E"; Value = "-1.01%" }
    @{ Row = 16; Col = "D"; Value = "0.04084" }
    @{ Row = 16; Col = "E"; Value = "-3.23%" }
    @{ Row = 17; Col = "D"; Value = "0.005822" }
    @{ Row = 17; Col = "E"; Value = "-3.12%" }
    @{ Row = 18; Col = "E"; Value = "-6.51%" }
    @{ Row = 19; Col = "D"; Value = "4.408" }
    @{ Row = 19; Col = "E"; Value = "2.56%" }
    @{ Row = 20; Col = "D"; Value = "0.3457" }
    @{ Row = 20; Col = "E"; Value = "-1.08%" }
    @{ Row = 21; Col = "D"; Value = "0.1384" }
    @{ Row = 21; Col = "E"; Value = "0.27%" }
    @{ Row = 22; Col = "D"; Value = "0.2582" }
    @{ Row = 22; Col = "E"; Value = "-2.81%" }
    @{ Row = 23; Col = "D"; Value = "0.001273" }
    @{ Row = 23; Col = "E"; Value = "1.70%" }
    @{ Row = 24; Col = "D"; Value = "0.004356" }
    @{ Row = 24; Col = "E"; Value = "-3.57%" }
    @{ Row = 25; Col = "D"; Value = "0.0001199" }
    @{ Row = 25; Col = "E"; Value = "-4.99%" }
    @{ Row = 26; Col = "D"; Value = "0.0003744" }
    @{ Row = 26; Col = "E"; Value = "-6.39%" }
    @{ Row = 38; Col = "D"; Value = "0.02590" }
    @{ Row = 38; Col = "E"; Value = "-0.48%" }
    @{ Row = 39; Col = "D"; Value = "0.05832" }
    @{ Row = 39; Col = "E"; Value = "8.45%" }
    @{ Row = 40; Col = "D"; Value = "0.007580" }
    @{ Row = 40; Col = "E"; Value = "-0.27%" }
    @{ Row = 41; Col = "E"; Value = "0.74%" }
    @{ Row = 42; Col = "D"; Value = "0.007160" }
    @{ Row = 42; Col = "E"; Value = "-1.86%" }
    @{ Row = 43; Col = "E"; Value = "-0.25%" }
    @{ Row = 44; Col = "D"; Value = "0.009078" }
    @{ Row = 44; Col = "E"; Value = "2.73%" }
    @{ Row = 45; Col = "D"; Value = "0.00006992" }
    @{ Row = 45; Col = "E"; Value = "0.84%" }
    @{ Row = 46; Col = "E"; Value = "-0.25%" }
    @{ Row = 47; Col = "D"; Value = "0.0005799" }
    @{ Row = 47; Col = "E"; Value = "-0.21%" }
    @{ Row = 48; Col = "D"; Value = "0.003530" }
    @{ Row = 48; Col = "E"; Value = "55.10%" }
    @{ Row = 49; Col = "D"; Value = "0.003537" }
    @{ Row = 49; Col = "E"; Value = "-2.70%" }
    @{ Row = 50; Col = "E"; Value = "-0.25%" }
    @{ Row = 51; Col = "E"; Value = "-0.25%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Col + $u.Row)
    # Force text storage first so Excel doesn't reinterpret the numeric-
    # looking literal (dropping trailing zeros / the trailing "%").
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
